$p = $ppt.ActivePresentation

# Remove the "Sara Demo" text from the title placeholder on slide 1,
# leaving an empty paragraph (matches the diff: run removed, endParaRPr kept).
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(1)
$titleLen = $titleShape.TextFrame.TextRange.Length
$titleShape.TextFrame.TextRange.Characters(1, $titleLen).Delete()

# Delete the second slide entirely (sldId 257 removed from sldIdLst,
# and its slide part dropped from the package).
$p.Slides.Item(2).Delete()
